$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_vals data (filtered save games) for rows 2-13.
# Column layout: A=date, B=TB, C=d2S, D=K, E=IP, F=Win, G=sum (=B+C+D+E)
$data = @(
    @{Row=2;  B=0.04271373187048222; C=3286.919754855326;  D=0.1494219747398047; E=10.19245300693656;  G=3297.304343568873}
    @{Row=3;  B=0.6606524410359556;  C=0.002571899574220771; D=0.1494219747398047; E=0.4942365360607697; G=1.306882851410751}
    @{Row=4;  B=3.286832544864788;   C=117.745847958593;   D=0.7527432677738641; E=10.19245300693656;  G=131.9778767781682}
    @{Row=5;  B=0.04271373187048222; C=0.306821227259698;  D=0.7527432677738641; E=0.4942365360607697; G=1.596514762964814}
    @{Row=6;  B=0.2917716402565462;  C=1.655778082260271;  D=261.3203778131603;  E=1133.036916526867;  G=1396.304844062544}
    @{Row=7;  B=3.286832544864788;   C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=8.974608811992548}
    @{Row=8;  B=1.455362044514542;   C=1.655778082260271;  D=0.1494219747398047; E=0.4942365360607697; G=3.754798637575387}
    @{Row=9;  B=0.6606524410359556;  C=1.655778082260271;  D=0.1494219747398047; E=0.4942365360607697; G=2.960089034096801}
    @{Row=10; B=3.286832544864788;   C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=8.974608811992548}
    @{Row=11; B=3.286832544864788;   C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=8.974608811992548}
    @{Row=12; B=0.6606524410359556;  C=1.655778082260271;  D=3.537761648806719;  E=0.4942365360607697; G=6.348428708163715}
    @{Row=13; B=1.455362044514542;   C=1.655778082260271;  D=0.7527432677738641; E=0.4942365360607697; G=4.358119930609447}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("G$r").Value = $entry.G
}
